# "started mu's in function" -- reformat the quantiles table:
#  - push everything down one row (adds a thin spacer row at the top)
#  - drop the old literal 0..5 index column (A) -- keep the column but blank it
#  - add two fresh blank rows at the bottom, matching the table's row style
#  - restyle the whole table: Cambria Math font, white fill, centered text,
#    a thick rule under the header / above the footer, thin rules around the
#    two "summary" rows (pre/post), taller rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Structural edits: insert a spacer row above the header, and append two
#    blank rows below the last data row.
# ---------------------------------------------------------------------------
$ws.Rows("1:1").Insert()

# Header now lives on row 2, data rows 3-8 (used to be 1-7 before the insert).
# Blank out the old manual index column (used to hold 0,1,2,3,4,5).
$ws.Range("A3:A8").ClearContents()

# Two new blank rows at the bottom of the table (rows 9 & 10).
$ws.Rows("9:10").Insert()

# ---------------------------------------------------------------------------
# 2. Fonts / fill / alignment shared by the whole table body (rows 1-10).
# ---------------------------------------------------------------------------
$body = $ws.Range("A1:J10")
$body.Font.Name = "Cambria Math"
$body.Font.Size = 11
$body.Interior.ThemeColor = 2
$body.Interior.TintAndShade = 0

$ws.Range("A2:J10").HorizontalAlignment = -4108
$ws.Range("A2:J10").VerticalAlignment = -4108

# Column B (the row labels) is bold throughout.
$ws.Range("B1:B10").Font.Bold = $true

# Header row (2) and footer row (8) are bold.
$ws.Range("B2:J2").Font.Bold = $true
$ws.Range("B8:J8").Font.Bold = $true

# ---------------------------------------------------------------------------
# 3. Borders.
# ---------------------------------------------------------------------------
# Thick rule under the spacer/header rows, and under the last data row.
$ws.Range("B1:J1").Borders.Item(9).LineStyle = 1
$ws.Range("B1:J1").Borders.Item(9).Weight = -4138
$ws.Range("A2:J2").Borders.Item(9).LineStyle = 1
$ws.Range("A2:J2").Borders.Item(9).Weight = -4138
$ws.Range("A8:J8").Borders.Item(9).LineStyle = 1
$ws.Range("A8:J8").Borders.Item(9).Weight = -4138

# Thin rules around the "pre" / "post" summary rows (4 & 5, and 7).
foreach ($r in 4, 5, 7) {
    $rng = $ws.Range("A" + $r + ":J" + $r)
    $rng.Borders.Item(8).LineStyle = 1
    $rng.Borders.Item(8).Weight = 2
    $rng.Borders.Item(9).LineStyle = 1
    $rng.Borders.Item(9).Weight = 2
}

# ---------------------------------------------------------------------------
# 4. Row heights / column widths.
# ---------------------------------------------------------------------------
$ws.Range("A1:J10").RowHeight = 27
$ws.Columns("A:A").ColumnWidth = 14.998697916666666
$ws.Columns("B:B").ColumnWidth = 14.998697916666666

# ---------------------------------------------------------------------------
# 5. Selection, matching the saved workbook's cursor position.
# ---------------------------------------------------------------------------
$ws.Range("C10").Select()

Write-Output "done"
